$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.688184333333333
$ws.Cells.Item(2, 8).Value = 5.064553
$ws.Cells.Item(2, 9).Value = 0.1174709603088462
$ws.Cells.Item(2, 10).Value = 0.1174709603088462
$ws.Cells.Item(2, 13).Value = 10.045207
$ws.Cells.Item(2, 14).Value = 30.135621
$ws.Cells.Item(2, 15).Value = 0.9365108453707793
$ws.Cells.Item(2, 16).Value = 0.9365108453707794
$ws.Cells.Item(2, 17).Value = 16.95816108249033
$ws.Cells.Item(2, 18).Value = 152.623449742413
$ws.Cells.Item(2, 19).Value = 0.1100128283453549
$ws.Cells.Item(2, 20).Value = 0.1100128283453549
$ws.Cells.Item(3, 7).Value = 1.688184333333333
$ws.Cells.Item(3, 8).Value = 5.064553
$ws.Cells.Item(3, 9).Value = 0.1174709603088462
$ws.Cells.Item(3, 10).Value = 0.1174709603088462
$ws.Cells.Item(3, 15).Value = 0.03971513502725754
$ws.Cells.Item(3, 16).Value = 0.03971513502725754
$ws.Cells.Item(3, 17).Value = 0.7191541459815556
$ws.Cells.Item(3, 18).Value = 6.472387313834
$ws.Cells.Item(3, 19).Value = 0.00466537505044744
$ws.Cells.Item(3, 20).Value = 0.00466537505044744
$ws.Cells.Item(4, 7).Value = 1.688184333333333
$ws.Cells.Item(4, 8).Value = 5.064553
$ws.Cells.Item(4, 9).Value = 0.1174709603088462
$ws.Cells.Item(4, 10).Value = 0.1174709603088462
$ws.Cells.Item(4, 13).Value = 0.255005
$ws.Cells.Item(4, 14).Value = 0.765015
$ws.Cells.Item(4, 15).Value = 0.02377401960196297
$ws.Cells.Item(4, 16).Value = 0.02377401960196297
$ws.Cells.Item(4, 17).Value = 0.4304954459216666
$ws.Cells.Item(4, 18).Value = 3.874459013295
$ws.Cells.Item(4, 19).Value = 0.002792756913043924
$ws.Cells.Item(4, 20).Value = 0.002792756913043924
$ws.Cells.Item(5, 9).Value = 0.07446899019451893
$ws.Cells.Item(5, 10).Value = 0.07446899019451893
$ws.Cells.Item(5, 13).Value = 10.045207
$ws.Cells.Item(5, 14).Value = 30.135621
$ws.Cells.Item(5, 15).Value = 0.9365108453707793
$ws.Cells.Item(5, 16).Value = 0.9365108453707794
$ws.Cells.Item(5, 17).Value = 10.75037718299767
$ws.Cells.Item(5, 18).Value = 96.753394646979
$ws.Cells.Item(5, 19).Value = 0.0697410169609772
$ws.Cells.Item(5, 20).Value = 0.0697410169609772
$ws.Cells.Item(6, 9).Value = 0.07446899019451893
$ws.Cells.Item(6, 10).Value = 0.07446899019451893
$ws.Cells.Item(6, 15).Value = 0.03971513502725754
$ws.Cells.Item(6, 16).Value = 0.03971513502725754
$ws.Cells.Item(6, 19).Value = 0.002957546000918837
$ws.Cells.Item(6, 20).Value = 0.002957546000918837
$ws.Cells.Item(7, 9).Value = 0.07446899019451893
$ws.Cells.Item(7, 10).Value = 0.07446899019451893
$ws.Cells.Item(7, 13).Value = 0.255005
$ws.Cells.Item(7, 14).Value = 0.765015
$ws.Cells.Item(7, 15).Value = 0.02377401960196297
$ws.Cells.Item(7, 16).Value = 0.02377401960196297
$ws.Cells.Item(7, 17).Value = 0.2729062659983333
$ws.Cells.Item(7, 18).Value = 2.456156393985
$ws.Cells.Item(7, 19).Value = 0.001770427232622881
$ws.Cells.Item(7, 20).Value = 0.001770427232622881
$ws.Cells.Item(8, 7).Value = 2.313935
$ws.Cells.Item(8, 8).Value = 6.941805
$ws.Cells.Item(8, 9).Value = 0.1610133213388724
$ws.Cells.Item(8, 10).Value = 0.1610133213388724
$ws.Cells.Item(8, 13).Value = 10.045207
$ws.Cells.Item(8, 14).Value = 30.135621
$ws.Cells.Item(8, 15).Value = 0.9365108453707793
$ws.Cells.Item(8, 16).Value = 0.9365108453707794
$ws.Cells.Item(8, 17).Value = 23.243956059545
$ws.Cells.Item(8, 18).Value = 209.195604535905
$ws.Cells.Item(8, 19).Value = 0.1507907216830243
$ws.Cells.Item(8, 20).Value = 0.1507907216830244
$ws.Cells.Item(9, 7).Value = 2.313935
$ws.Cells.Item(9, 8).Value = 6.941805
$ws.Cells.Item(9, 9).Value = 0.1610133213388724
$ws.Cells.Item(9, 10).Value = 0.1610133213388724
$ws.Cells.Item(9, 15).Value = 0.03971513502725754
$ws.Cells.Item(9, 16).Value = 0.03971513502725754
$ws.Cells.Item(9, 17).Value = 0.9857193411433333
$ws.Cells.Item(9, 18).Value = 8.87147407029
$ws.Cells.Item(9, 19).Value = 0.006394665798160526
$ws.Cells.Item(9, 20).Value = 0.006394665798160526
$ws.Cells.Item(10, 7).Value = 2.313935
$ws.Cells.Item(10, 8).Value = 6.941805
$ws.Cells.Item(10, 9).Value = 0.1610133213388724
$ws.Cells.Item(10, 10).Value = 0.1610133213388724
$ws.Cells.Item(10, 13).Value = 0.255005
$ws.Cells.Item(10, 14).Value = 0.765015
$ws.Cells.Item(10, 15).Value = 0.02377401960196297
$ws.Cells.Item(10, 16).Value = 0.02377401960196297
$ws.Cells.Item(10, 17).Value = 0.590064994675
$ws.Cells.Item(10, 18).Value = 5.310584952075
$ws.Cells.Item(10, 19).Value = 0.003827933857687515
$ws.Cells.Item(10, 20).Value = 0.003827933857687515
$ws.Cells.Item(11, 7).Value = 0.5887196666666666
$ws.Cells.Item(11, 8).Value = 1.766159
$ws.Cells.Item(11, 9).Value = 0.04096558843161708
$ws.Cells.Item(11, 10).Value = 0.04096558843161708
$ws.Cells.Item(11, 13).Value = 10.045207
$ws.Cells.Item(11, 14).Value = 30.135621
$ws.Cells.Item(11, 15).Value = 0.9365108453707793
$ws.Cells.Item(11, 16).Value = 0.9365108453707794
$ws.Cells.Item(11, 17).Value = 5.913810916637666
$ws.Cells.Item(11, 18).Value = 53.224298249739
$ws.Cells.Item(11, 19).Value = 0.03836471785320513
$ws.Cells.Item(11, 20).Value = 0.03836471785320513
$ws.Cells.Item(12, 7).Value = 0.5887196666666666
$ws.Cells.Item(12, 8).Value = 1.766159
$ws.Cells.Item(12, 9).Value = 0.04096558843161708
$ws.Cells.Item(12, 10).Value = 0.04096558843161708
$ws.Cells.Item(12, 15).Value = 0.03971513502725754
$ws.Cells.Item(12, 16).Value = 0.03971513502725754
$ws.Cells.Item(12, 17).Value = 0.2507902607224444
$ws.Cells.Item(12, 18).Value = 2.257112346502
$ws.Cells.Item(12, 19).Value = 0.001626953876032732
$ws.Cells.Item(12, 20).Value = 0.001626953876032732
$ws.Cells.Item(13, 7).Value = 0.5887196666666666
$ws.Cells.Item(13, 8).Value = 1.766159
$ws.Cells.Item(13, 9).Value = 0.04096558843161708
$ws.Cells.Item(13, 10).Value = 0.04096558843161708
$ws.Cells.Item(13, 13).Value = 0.255005
$ws.Cells.Item(13, 14).Value = 0.765015
$ws.Cells.Item(13, 15).Value = 0.02377401960196297
$ws.Cells.Item(13, 16).Value = 0.02377401960196297
$ws.Cells.Item(13, 17).Value = 0.1501264585983333
$ws.Cells.Item(13, 18).Value = 1.351138127385
$ws.Cells.Item(13, 19).Value = 0.0009739167023792118
$ws.Cells.Item(13, 20).Value = 0.0009739167023792118
$ws.Cells.Item(14, 7).Value = 7.350912666666666
$ws.Cells.Item(14, 8).Value = 22.052738
$ws.Cells.Item(14, 9).Value = 0.5115073946899924
$ws.Cells.Item(14, 10).Value = 0.5115073946899924
$ws.Cells.Item(14, 13).Value = 10.045207
$ws.Cells.Item(14, 14).Value = 30.135621
$ws.Cells.Item(14, 15).Value = 0.9365108453707793
$ws.Cells.Item(14, 16).Value = 0.9365108453707794
$ws.Cells.Item(14, 17).Value = 73.84143937558866
$ws.Cells.Item(14, 18).Value = 664.5729543802979
$ws.Cells.Item(14, 19).Value = 0.4790322226145297
$ws.Cells.Item(14, 20).Value = 0.4790322226145297
$ws.Cells.Item(15, 7).Value = 7.350912666666666
$ws.Cells.Item(15, 8).Value = 22.052738
$ws.Cells.Item(15, 9).Value = 0.5115073946899924
$ws.Cells.Item(15, 10).Value = 0.5115073946899924
$ws.Cells.Item(15, 15).Value = 0.03971513502725754
$ws.Cells.Item(15, 16).Value = 0.03971513502725754
$ws.Cells.Item(15, 17).Value = 3.131434889307111
$ws.Cells.Item(15, 18).Value = 28.182914003764
$ws.Cells.Item(15, 19).Value = 0.02031458524755377
$ws.Cells.Item(15, 20).Value = 0.02031458524755377
$ws.Cells.Item(16, 7).Value = 7.350912666666666
$ws.Cells.Item(16, 8).Value = 22.052738
$ws.Cells.Item(16, 9).Value = 0.5115073946899924
$ws.Cells.Item(16, 10).Value = 0.5115073946899924
$ws.Cells.Item(16, 13).Value = 0.255005
$ws.Cells.Item(16, 14).Value = 0.765015
$ws.Cells.Item(16, 15).Value = 0.02377401960196297
$ws.Cells.Item(16, 16).Value = 0.02377401960196297
$ws.Cells.Item(16, 17).Value = 1.874519484563333
$ws.Cells.Item(16, 18).Value = 16.87067536107
$ws.Cells.Item(16, 19).Value = 0.01216058682790889
$ws.Cells.Item(16, 20).Value = 0.01216058682790889
$ws.Cells.Item(17, 7).Value = 1.359126666666667
$ws.Cells.Item(17, 8).Value = 4.07738
$ws.Cells.Item(17, 9).Value = 0.09457374503615294
$ws.Cells.Item(17, 10).Value = 0.09457374503615293
$ws.Cells.Item(17, 13).Value = 10.045207
$ws.Cells.Item(17, 14).Value = 30.135621
$ws.Cells.Item(17, 15).Value = 0.9365108453707793
$ws.Cells.Item(17, 16).Value = 0.9365108453707794
$ws.Cells.Item(17, 17).Value = 13.65270870588667
$ws.Cells.Item(17, 18).Value = 122.87437835298
$ws.Cells.Item(17, 19).Value = 0.08856933791368814
$ws.Cells.Item(17, 20).Value = 0.08856933791368814
$ws.Cells.Item(18, 7).Value = 1.359126666666667
$ws.Cells.Item(18, 8).Value = 4.07738
$ws.Cells.Item(18, 9).Value = 0.09457374503615294
$ws.Cells.Item(18, 10).Value = 0.09457374503615293
$ws.Cells.Item(18, 15).Value = 0.03971513502725754
$ws.Cells.Item(18, 16).Value = 0.03971513502725754
$ws.Cells.Item(18, 17).Value = 0.5789779930711111
$ws.Cells.Item(18, 18).Value = 5.21080193764
$ws.Cells.Item(18, 19).Value = 0.003756009054144242
$ws.Cells.Item(18, 20).Value = 0.003756009054144241
$ws.Cells.Item(19, 7).Value = 1.359126666666667
$ws.Cells.Item(19, 8).Value = 4.07738
$ws.Cells.Item(19, 9).Value = 0.09457374503615294
$ws.Cells.Item(19, 10).Value = 0.09457374503615293
$ws.Cells.Item(19, 13).Value = 0.255005
$ws.Cells.Item(19, 14).Value = 0.765015
$ws.Cells.Item(19, 15).Value = 0.02377401960196297
$ws.Cells.Item(19, 16).Value = 0.02377401960196297
$ws.Cells.Item(19, 17).Value = 0.3465840956333333
$ws.Cells.Item(19, 18).Value = 3.1192568607
$ws.Cells.Item(19, 19).Value = 0.002248398068320548
$ws.Cells.Item(19, 20).Value = 0.002248398068320548
